# Auto-generated script to update Leve profit calculation cells
# per scheduled market-data refresh (see commit message).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 982
$ws.Range("J19").Value = 922.5
$ws.Range("L19").Value = 922.5
$ws.Range("N19").Value = -1272.5
$ws.Range("H113").Value = 5721.636
$ws.Range("I113").Value = 4648
$ws.Range("J113").Value = 6124.25
$ws.Range("K113").Value = 4648
$ws.Range("L113").Value = 6124.25
$ws.Range("M113").Value = -1394
$ws.Range("N113").Value = -12632.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 15964.333
$ws.Range("J44").Value = 15964.333
$ws.Range("L44").Value = 15964.333
$ws.Range("N44").Value = -16940.333
$ws.Range("H45").Value = 10151.286
$ws.Range("I45").Value = 10711.583
$ws.Range("J45").Value = 6789.5
$ws.Range("K45").Value = 10711.583
$ws.Range("L45").Value = 6789.5
$ws.Range("M45").Value = -10334.583
$ws.Range("N45").Value = -7543.5
$ws.Range("H63").Value = 2604.2104
$ws.Range("I63").Value = 1822
$ws.Range("K63").Value = 1822
$ws.Range("M63").Value = -1136
$ws.Range("H66").Value = 2604.2104
$ws.Range("I66").Value = 1822
$ws.Range("K66").Value = 9110
$ws.Range("M66").Value = -5678
$ws.Range("H88").Value = 1421.4117
$ws.Range("I88").Value = 871.7
$ws.Range("J88").Value = 2206.7144
$ws.Range("K88").Value = 871.7
$ws.Range("L88").Value = 2206.7144
$ws.Range("M88").Value = -465.7
$ws.Range("N88").Value = -3018.7144
$ws.Range("H91").Value = 1421.4117
$ws.Range("I91").Value = 871.7
$ws.Range("J91").Value = 2206.7144
$ws.Range("K91").Value = 871.7
$ws.Range("L91").Value = 2206.7144
$ws.Range("M91").Value = 532.3
$ws.Range("N91").Value = -5014.7144
$ws.Range("H141").Value = 2465.2942
$ws.Range("I141").Value = 2425.625
$ws.Range("K141").Value = 7276.875
$ws.Range("M141").Value = -2096.875

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H99").Value = 2747.0322
$ws.Range("I99").Value = 2778.8462
$ws.Range("J99").Value = 2581.6
$ws.Range("K99").Value = 2778.8462
$ws.Range("L99").Value = 2581.6
$ws.Range("M99").Value = -1280.8462
$ws.Range("N99").Value = -5577.6
$ws.Range("H105").Value = 3379.6086
$ws.Range("I105").Value = 3551.7222
$ws.Range("K105").Value = 3551.7222
$ws.Range("M105").Value = -1804.7222
$ws.Range("H132").Value = 5811.5557
$ws.Range("I132").Value = 5907.2954
$ws.Range("K132").Value = 17721.8862
$ws.Range("M132").Value = -15191.8862

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3272.5312
$ws.Range("I16").Value = 3433.1365
$ws.Range("J16").Value = 2919.2
$ws.Range("K16").Value = 3433.1365
$ws.Range("L16").Value = 2919.2
$ws.Range("M16").Value = -3146.1365
$ws.Range("N16").Value = -3493.2
$ws.Range("H31").Value = 3272.9583
$ws.Range("I31").Value = 2629.9375
$ws.Range("K31").Value = 2629.9375
$ws.Range("M31").Value = -2334.9375
$ws.Range("H34").Value = 3272.9583
$ws.Range("I34").Value = 2629.9375
$ws.Range("K34").Value = 2629.9375
$ws.Range("M34").Value = -2427.9375
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H97").Value = 68560
$ws.Range("J97").Value = 67464.836
$ws.Range("L97").Value = 67464.836
$ws.Range("N97").Value = -69446.836
$ws.Range("H99").Value = 11064.393
$ws.Range("J99").Value = 11725.272
$ws.Range("L99").Value = 11725.272
$ws.Range("N99").Value = -14721.272
$ws.Range("H113").Value = 3272.5312
$ws.Range("I113").Value = 3433.1365
$ws.Range("J113").Value = 2919.2
$ws.Range("K113").Value = 3433.1365
$ws.Range("L113").Value = 2919.2
$ws.Range("M113").Value = -1263.1365
$ws.Range("N113").Value = -7259.2
$ws.Range("H122").Value = 2320.1
$ws.Range("J122").Value = 2498
$ws.Range("L122").Value = 7494
$ws.Range("N122").Value = -12394
$ws.Range("H126").Value = 11064.393
$ws.Range("J126").Value = 11725.272
$ws.Range("L126").Value = 35175.81600000001
$ws.Range("N126").Value = -40115.81600000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2042.697
$ws.Range("J5").Value = 2526.9092
$ws.Range("L5").Value = 7580.7276
$ws.Range("N5").Value = -7804.7276
$ws.Range("H25").Value = 6516.75
$ws.Range("I25").Value = 8925.125
$ws.Range("J25").Value = 1700
$ws.Range("K25").Value = 26775.375
$ws.Range("L25").Value = 5100
$ws.Range("M25").Value = -26606.375
$ws.Range("N25").Value = -5438
$ws.Range("H30").Value = 6516.75
$ws.Range("I30").Value = 8925.125
$ws.Range("J30").Value = 1700
$ws.Range("K30").Value = 26775.375
$ws.Range("L30").Value = 5100
$ws.Range("M30").Value = -26673.375
$ws.Range("N30").Value = -5304
$ws.Range("H32").Value = 1662.5
$ws.Range("J32").Value = 2325
$ws.Range("L32").Value = 6975
$ws.Range("N32").Value = -7541
$ws.Range("H55").Value = 2934
$ws.Range("J55").Value = 2934
$ws.Range("L55").Value = 8802
$ws.Range("N55").Value = -9156
$ws.Range("H121").Value = 1973.6957
$ws.Range("I121").Value = 1542.0526
$ws.Range("J121").Value = 4024
$ws.Range("K121").Value = 4626.1578
$ws.Range("L121").Value = 12072
$ws.Range("M121").Value = -3316.1578
$ws.Range("N121").Value = -14692

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 312891.75
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 312891.75
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 312891.75
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -313169.75
$ws.Range("H126").Value = 5753.1665
$ws.Range("I126").Value = 5440.3335
$ws.Range("K126").Value = 16321.0005
$ws.Range("M126").Value = -13851.0005
$ws.Range("H135").Value = 2042.697
$ws.Range("J135").Value = 2526.9092
$ws.Range("L135").Value = 22742.1828
$ws.Range("N135").Value = -27812.1828

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 100000
$ws.Range("J20").Value = 100000
$ws.Range("L20").Value = 100000
$ws.Range("N20").Value = -100452
$ws.Range("H22").Value = 2764
$ws.Range("I22").Value = 2481.75
$ws.Range("K22").Value = 2481.75
$ws.Range("M22").Value = -2186.75
$ws.Range("H27").Value = 2764
$ws.Range("I27").Value = 2481.75
$ws.Range("K27").Value = 2481.75
$ws.Range("M27").Value = -2374.75
$ws.Range("H68").Value = 2817.2856
$ws.Range("I68").Value = 2534.1155
$ws.Range("J68").Value = 6498.5
$ws.Range("K68").Value = 2534.1155
$ws.Range("L68").Value = 6498.5
$ws.Range("M68").Value = -1785.1155
$ws.Range("N68").Value = -7996.5
$ws.Range("H71").Value = 2817.2856
$ws.Range("I71").Value = 2534.1155
$ws.Range("J71").Value = 6498.5
$ws.Range("K71").Value = 12670.5775
$ws.Range("L71").Value = 32492.5
$ws.Range("M71").Value = -8926.577499999999
$ws.Range("N71").Value = -39980.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 18249.25
$ws.Range("J14").Value = 17499
$ws.Range("L14").Value = 17499
$ws.Range("N14").Value = -17835
$ws.Range("H62").Value = 5749.75
$ws.Range("I62").Value = 5999.5
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 5999.5
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -5375.5
$ws.Range("N62").Value = -6748
$ws.Range("H65").Value = 5749.75
$ws.Range("I65").Value = 5999.5
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 29997.5
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -26877.5
$ws.Range("N65").Value = -33740
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H96").Value = 4134.0713
$ws.Range("I96").Value = 4679.727
$ws.Range("J96").Value = 2133.3333
$ws.Range("K96").Value = 4679.727
$ws.Range("L96").Value = 2133.3333
$ws.Range("M96").Value = -3306.727
$ws.Range("N96").Value = -4879.3333
$ws.Range("H122").Value = 4310.1904
$ws.Range("I122").Value = 2949.2942
$ws.Range("K122").Value = 8847.882599999999
$ws.Range("M122").Value = -6397.882599999999
